$d = $word.ActiveDocument

function Esc([string]$s) {
    return $s.Replace("&","&amp;").Replace("<","&lt;").Replace(">","&gt;")
}

$pkgOpen  = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Paragraph 3: "[team#]" -> "[teamName]" (spellStart/spellEnd added around
# the word, gramStart/gramEnd kept, trailing "#" dropped from the last run)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$xml3 = $pkgOpen +
  '<w:p w:rsidR="00351E4D" w:rsidRDefault="00BB015D" w:rsidP="003C4C89">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>teamName</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>]</w:t></w:r>' +
  '</w:p>' + $pkgClose
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Paragraph 4: "[section#]" -> "Section 1" (drop the brackets/proofErr,
# relocate the _GoBack bookmark here, between "Section " and "1")
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$xml4 = $pkgOpen +
  '<w:p w:rsidR="00BB015D" w:rsidRDefault="00BB015D" w:rsidP="003C4C89">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Section </w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>1</w:t></w:r>' +
  '</w:p>' + $pkgClose
$p4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# Paragraph 5: "2" -> "2" + "hrs" (new trailing run)
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$xml5 = $pkgOpen +
  '<w:p w:rsidR="00BB015D" w:rsidRDefault="00BB015D" w:rsidP="003C4C89">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>2</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>hrs</w:t></w:r>' +
  '</w:p>' + $pkgClose
$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------------
# Last paragraph: drop the now-relocated _GoBack bookmark (id 0) from the
# end of the weekly-summary paragraph; the run text itself is unchanged, so
# pull it back from the live document to avoid retyping the long passage.
# ---------------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$fullText = $pLast.Range.Text
$fullText = $fullText.TrimEnd([char]13, [char]7)
$marker = "Wednesdays at "
$splitAt = $fullText.IndexOf($marker) + $marker.Length
$part1 = $fullText.Substring(0, $splitAt)
$part2 = $fullText.Substring($splitAt)

$xmlLast = $pkgOpen +
  '<w:p w:rsidR="00BB015D" w:rsidRPr="003C4C89" w:rsidRDefault="00BB015D" w:rsidP="003C4C89">' +
    '<w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">' + (Esc $part1) + '</w:t></w:r>' +
    '<w:r w:rsidR="00A97B7A"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>' + (Esc $part2) + '</w:t></w:r>' +
  '</w:p>' + $pkgClose
$pLast.Range.InsertXML($xmlLast)
